$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 2.4
$ws.Range("H2").Value = 2.75
$ws.Range("I2").Value = 3.5
$ws.Range("K2").Value = 1.8
$ws.Range("U2").Value = 2.5
$ws.Range("V2").Value = 1.5
$ws.Range("X2").Value = 9.5
$ws.Range("AA2").Value = 29
$ws.Range("AB2").Value = 51
$ws.Range("AD2").Value = 6
$ws.Range("AF2").Value = 101
$ws.Range("AJ2").Value = 15
$ws.Range("AL2").Value = 41
$ws.Range("AO2").Value = 17
$ws.Range("AT2").Value = 2
$ws.Range("AV2").Value = 101
$ws.Range("AX2").Value = 23
$ws.Range("BA2").Value = 151
